# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps that were refreshed when the
# handback report was (re)generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# These are stored as plain text strings (not real Excel dates), so force
# text assignment with a leading-quote-safe approach by setting .Value on a
# cell that's already formatted as text/general, writing the literal string.

$wsOverview.Range("G3").Value = "2016-08-27 14:44:57"

$wsZhCn.Range("H3").Value = "2016-08-27 14:44:53"
$wsZhCn.Range("K3").Value = "2016-08-27 14:45:16"

$wsDeDe.Range("K3").Value = "2016-08-27 14:45:23"
